# The "depth" example sheets (Biosample, Biosample1) previously had a single
# header cell "depth" in A1. The commit wraps the example in an outer
# "biosamples:" key, which means the attribute row now needs an "id" column
# before "depth", i.e. A1 = "id", B1 = "depth".
$wb = $excel.ActiveWorkbook

$sheetNames = @("Biosample", "Biosample1")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "depth"
    $ws.Range("A1").Value = "id"
}
